$d = $word.ActiveDocument

# Locate the exact point between '...see "you lose" text' and the
# trailing '.' run -- this is where the _GoBack bookmark currently sits
# (between the "text" run and the "." run), in the
# '- I should see "you lose" text.' bullet under the second
# "you lose" Acceptance Criteria.
$quote1 = [char]8220
$quote2 = [char]8221
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "see " + $quote1 + "you lose" + $quote2 + " text", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'see <<you lose>> text' run"
}
$splitPos = $findRange.End

# Find the 1-based paragraph index that contains $splitPos.
$origIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Start -le $splitPos -and $splitPos -le $p.Range.End) {
        $origIndex = $idx
    }
}
if ($origIndex -eq 0) {
    throw "Could not find paragraph containing split position"
}

# Step 1: split the paragraph right at the bookmark's position. Word
# clones the current (numbered list-item) paragraph formatting onto the
# brand-new paragraph, so the bookmark (and the trailing ".") move into
# that new paragraph while the original paragraph keeps the list-item
# formatting.
$breakRange = $d.Range($splitPos, $splitPos)
$breakRange.InsertParagraphAfter()

$newIndex = $origIndex + 1
$origPara = $d.Paragraphs.Item($origIndex)
$newPara = $d.Paragraphs.Item($newIndex)

# Step 2: the new paragraph inherited the numbered list-item formatting;
# fix it up to the plain "ListParagraph / ind left=1080" formatting the
# other blank paragraphs in this list use.
$newPara.Range.ListFormat.RemoveNumbers()
$newPara.Style = "List Paragraph"
$newPara.Range.ParagraphFormat.LeftIndent = 54

# Step 3: move the "." back to the end of the original paragraph (before
# its paragraph mark), leaving only the bookmark behind in the new
# (now blank) paragraph.
$newPara = $d.Paragraphs.Item($newIndex)
$dotRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + 1)
$dotRange.Delete()
$origPara = $d.Paragraphs.Item($origIndex)
$insertAt = $origPara.Range.End - 1
$d.Range($insertAt, $insertAt).InsertAfter(".")

# Step 4: two blank "ListParagraph" paragraphs used to follow -- now that
# the bookmark lives in its own new blank paragraph, those two old blank
# paragraphs are redundant leftovers from the old use-case diagram and
# should both be removed.
$extra1 = $d.Paragraphs.Item($newIndex + 1)
$extra1.Range.Delete()
$extra2 = $d.Paragraphs.Item($newIndex + 1)
$extra2.Range.Delete()
